# Generate Report for Handoff
# Adds two new "Ready for handoff" rows (163d6d51-... and e72297d2-...)
# to the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$guid1 = "163d6d51-a136-4906-a1a5-d29efb81317b"
$guid2 = "e72297d2-071a-4b71-bb90-a4ed5bad548c"
$hash1 = "3355fda8460681208c51f640e37d13cba5f7fde2"
$hash2 = "ca4216852a7f8593b4c5f6005e789eb3a4fc8779"

$status = "Ready for handoff"
$ext = ".md"
$zeroDate = "0001-01-01 00:00:00"
$reason = "Include"

$zhDate = "2016-03-14 08:57:16"
$deDate = "2016-03-14 08:57:24"
$overviewDate = "2016-57-14 08:57:24"

$md1 = $guid1 + ".md"
$md2 = $guid2 + ".md"

$xlfZh1 = $guid1 + "." + $hash1 + ".zh-cn.xlf"
$xlfZh2 = $guid2 + "." + $hash2 + ".zh-cn.xlf"
$xlfDe1 = $guid1 + "." + $hash1 + ".de-de.xlf"
$xlfDe2 = $guid2 + "." + $hash2 + ".de-de.xlf"

$orgRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/9a74638d1cebff8e110b6258c1af874ca78c4d48/e2e/"
$handoffZhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e8bee5e51bbc699f92d584505897b8a0fc210577/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/"
$handoffDeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c506d818b04de5fc71a82761e1a8cf777096d0c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $md1
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $overviewDate

$wsOverview.Range("A5").Value = $md2
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $orgRepoBase + $md1, [Type]::Missing, [Type]::Missing, $md1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $orgRepoBase + $md2, [Type]::Missing, [Type]::Missing, $md2)

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File
# | Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $md1
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = $xlfZh1
$wsZh.Range("E4").Value = $zhDate
$wsZh.Range("H4").Value = $zeroDate
$wsZh.Range("I4").Value = $reason

$wsZh.Range("A5").Value = $md2
$wsZh.Range("B5").Value = $ext
$wsZh.Range("C5").Value = $status
$wsZh.Range("D5").Value = $xlfZh2
$wsZh.Range("E5").Value = $zhDate
$wsZh.Range("H5").Value = $zeroDate
$wsZh.Range("I5").Value = $reason

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $orgRepoBase + $md1, [Type]::Missing, [Type]::Missing, $md1)
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $orgRepoBase + $md1, [Type]::Missing, [Type]::Missing, $ext)
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $handoffZhBase + $xlfZh1, [Type]::Missing, [Type]::Missing, $xlfZh1)

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $orgRepoBase + $md2, [Type]::Missing, [Type]::Missing, $md2)
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), $orgRepoBase + $md2, [Type]::Missing, [Type]::Missing, $ext)
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), $handoffZhBase + $xlfZh2, [Type]::Missing, [Type]::Missing, $xlfZh2)

# ---------------------------------------------------------------------------
# de-de sheet: same column layout as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $md1
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = $xlfDe1
$wsDe.Range("E4").Value = $deDate
$wsDe.Range("H4").Value = $zeroDate
$wsDe.Range("I4").Value = $reason

$wsDe.Range("A5").Value = $md2
$wsDe.Range("B5").Value = $ext
$wsDe.Range("C5").Value = $status
$wsDe.Range("D5").Value = $xlfDe2
$wsDe.Range("E5").Value = $deDate
$wsDe.Range("H5").Value = $zeroDate
$wsDe.Range("I5").Value = $reason

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $orgRepoBase + $md1, [Type]::Missing, [Type]::Missing, $md1)
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $orgRepoBase + $md1, [Type]::Missing, [Type]::Missing, $ext)
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $handoffDeBase + $xlfDe1, [Type]::Missing, [Type]::Missing, $xlfDe1)

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $orgRepoBase + $md2, [Type]::Missing, [Type]::Missing, $md2)
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), $orgRepoBase + $md2, [Type]::Missing, [Type]::Missing, $ext)
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), $handoffDeBase + $xlfDe2, [Type]::Missing, [Type]::Missing, $xlfDe2)

"done"
